# T460 - update trip log (foaie de parcurs) values for June 2022

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Starting km counter
$ws.Range("B12").Value = 356153

# Daily trip rows (A=day, B=km, C=place, D=observation)
$ws.Range("B15").Value = 30
$ws.Range("C15").Value = "Acasa-Birou"
$ws.Range("D15").Value = "Interes Serviciu"

$ws.Range("B16").Value = 85
$ws.Range("C16").Value = "Cluj-Apahida"
$ws.Range("D16").Value = "Interes Serviciu"

$ws.Range("B19").Value = 85
$ws.Range("C19").Value = "Cluj-Apahida"
$ws.Range("D19").Value = "Interes Serviciu"

$ws.Range("B20").Value = 30
$ws.Range("C20").Value = "Acasa-Birou"
$ws.Range("D20").Value = " "

$ws.Range("B21").Value = 30
$ws.Range("C21").Value = "Acasa-Birou"
$ws.Range("D21").Value = " "

$ws.Range("B22").Value = 421
$ws.Range("C22").Value = "Cluj-Satu-Mare"
$ws.Range("D22").Value = "Interes Serviciu"

$ws.Range("B23").Value = 85
$ws.Range("C23").Value = "Cluj-Apahida"
$ws.Range("D23").Value = "Interes Serviciu"

$ws.Range("B27").Value = 30
$ws.Range("C27").Value = "Acasa-Birou"
$ws.Range("D27").Value = " "

$ws.Range("B28").Value = 30
$ws.Range("C28").Value = "Acasa-Birou"
$ws.Range("D28").Value = " "

$ws.Range("B29").Value = 356
$ws.Range("C29").Value = "Cluj-Baia-Mare"
$ws.Range("D29").Value = "Interes Serviciu"

$ws.Range("B30").Value = 257
$ws.Range("C30").Value = "Cluj-Bistrita"
$ws.Range("D30").Value = "Interes Serviciu"

$ws.Range("B33").Value = 156
$ws.Range("C33").Value = "Cluj-Zalau"
$ws.Range("D33").Value = "Interes Serviciu"

$ws.Range("B34").Value = 121
$ws.Range("C34").Value = "Cluj-Turda"
$ws.Range("D34").Value = "Interes Serviciu"

$ws.Range("B35").Value = 30
$ws.Range("C35").Value = "Acasa-Birou"
$ws.Range("D35").Value = " "

$ws.Range("B36").Value = 30
$ws.Range("C36").Value = "Acasa-Birou"
$ws.Range("D36").Value = " "

$ws.Range("B37").Value = 421
$ws.Range("C37").Value = "Cluj-Satu-Mare"
$ws.Range("D37").Value = "Interes Serviciu"

$ws.Range("B40").Value = 121
$ws.Range("C40").Value = "Cluj-Turda"
$ws.Range("D40").Value = "Interes Serviciu"

$ws.Range("B41").Value = 421
$ws.Range("C41").Value = "Cluj-Satu-Mare"
$ws.Range("D41").Value = "Interes Serviciu"

$ws.Range("B42").Value = 101
$ws.Range("C42").Value = "Cluj-Dej"
$ws.Range("D42").Value = "Interes Serviciu"

$ws.Range("B43").Value = 152
$ws.Range("C43").Value = "Cluj-Cmp. Turzii"
$ws.Range("D43").Value = "Interes Serviciu"

# Totals
$ws.Range("B44").Value = 2992
$ws.Range("B45").Value = 359145
